$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.633.56'
$ws.Range("E2").Value = '  +1.07%  '

$ws.Range("D3").Value = '3.422.19'
$ws.Range("E3").Value = '  -0.13%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.89'
$ws.Range("E5").Value = '  +0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.73'
$ws.Range("E6").Value = '  +4.61%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("E7").Value = '  +5.68%  '

$ws.Range("D9").Value = '3.419.23'
$ws.Range("E9").Value = '  -0.11%  '

$ws.Range("E10").Value = '  +1.79%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.97'
$ws.Range("E11").Value = '  +1.89%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.414'
$ws.Range("E12").Value = '  +1.19%  '

$ws.Range("D13").Value = '4.016.48'
$ws.Range("E13").Value = '  -0.07%  '

$ws.Range("E14").Value = '  +0.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.25'
$ws.Range("E15").Value = '  -2.05%  '

$ws.Range("D16").Value = '66.579.07'
$ws.Range("E16").Value = '  +0.92%  '

$ws.Range("E17").Value = '  +1.80%  '

$ws.Range("D18").Value = '3.430.60'
$ws.Range("E18").Value = '  +0.34%  '

$ws.Range("E19").Value = '  +0.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.85'
$ws.Range("E20").Value = '  +1.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '368.72'
$ws.Range("E21").Value = '  +1.03%  '

$ws.Range("E22").Value = '  -1.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.19'
$ws.Range("E23").Value = '  +2.74%  '

$ws.Range("E24").Value = '  +0.45%  '

$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000126'
$ws.Range("E25").Value = '  +6.66%  '

$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.534'
$ws.Range("E26").Value = '  +1.81%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.85'
$ws.Range("E27").Value = '  +1.80%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.179'
$ws.Range("E28").Value = '  +1.47%  '

$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.80'
$ws.Range("E30").Value = '  +0.29%  '

$ws.Range("E31").Value = '  +0.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.35'
$ws.Range("E32").Value = '  -1.97%  '

$ws.Range("E33").Value = '  +0.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.06'
$ws.Range("E34").Value = '  +0.38%  '

$ws.Range("E35").Value = '  -1.94%  '

$ws.Range("E36").Value = '  +0.67%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.63'
$ws.Range("E37").Value = '  +2.43%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.867'
$ws.Range("E38").Value = '  -1.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.54'
$ws.Range("E39").Value = '  -5.11%  '

$ws.Range("E40").Value = '  +2.13%  '

$ws.Range("E41").Value = '  +4.64%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.42'
$ws.Range("E42").Value = '  +0.66%  '

$ws.Range("D43").Value = '2.708.64'
$ws.Range("E43").Value = '  +0.32%  '

$ws.Range("E44").Value = '  -0.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0689'
$ws.Range("E45").Value = '  +1.36%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.03'
$ws.Range("E46").Value = '  +4.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '337.87'
$ws.Range("E47").Value = '  +10.67%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '39.90'
$ws.Range("E48").Value = '  -0.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0286'
$ws.Range("E49").Value = '  -0.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.104'
$ws.Range("E50").Value = '  +3.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '32.10'
$ws.Range("E51").Value = '  +6.12%  '
